# Fix up the "汽車" (car/vehicle) property sheet so it has the same
# column layout (headers + metadata columns) as every other sheet in
# this legislator's property-disclosure workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# ---- Header row (row 1): proper column labels instead of duplicated data ----
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ---- Data row (row 2): keep existing values, append the metadata columns ----
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2012-04-26"
$ws.Range("K2").Value = "黃昭順"
$ws.Range("L2").Value = 665
$ws.Range("M2").Value = "tmp4c4f1"
$ws.Range("N2").Value = 29

# Match the header/data styling already used by columns B-G on this sheet
$ws.Range("H1:N1").Style = $ws.Range("B1").Style
$ws.Range("H2:N2").Style = $ws.Range("B2").Style
